$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.370.86"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "2.652.14"
$ws.Range("E3").Value = "  +2.89%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.13"
$ws.Range("E5").Value = "  +3.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.71"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").Value = "2.651.43"
$ws.Range("E9").Value = "  +2.84%  "
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.65"
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("E13").Value = "  +3.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.39"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").Value = "3.126.07"
$ws.Range("E15").Value = "  +2.69%  "
$ws.Range("D16").Value = "63.200.21"
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("D18").Value = "2.646.49"
$ws.Range("E18").Value = "  +2.43%  "
$ws.Range("E19").Value = "  +2.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.45"
$ws.Range("E20").Value = "  +3.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "343.41"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.85"
$ws.Range("E22").Value = "  +3.18%  "
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.19"
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("E25").Value = "  +1.64%  "
$ws.Range("E26").Value = "  -2.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.68"
$ws.Range("E27").Value = "  +5.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.164"
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "549.24"
$ws.Range("E29").Value = "  +16.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.90"
$ws.Range("E31").Value = "  -0.58%  "
$ws.Range("E32").Value = "  +3.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.80"
$ws.Range("E33").Value = "  +6.61%  "
$ws.Range("D34").Value = "0.0₃0812"
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "172.51"
$ws.Range("E35").Value = "  -2.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.10"
$ws.Range("E36").Value = "  +12.02%  "
$ws.Range("E37").Value = "  +1.85%  "
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.14"
$ws.Range("E39").Value = "  +1.43%  "
$ws.Range("E40").Value = "  +6.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "172.63"
$ws.Range("E41").Value = "  +7.68%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.76"
$ws.Range("E43").Value = "  +1.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.43"
$ws.Range("E44").Value = "  +3.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0572"
$ws.Range("E45").Value = "  +6.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.633"
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("E47").Value = "  +1.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0962"
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.82"
$ws.Range("E49").Value = "  +3.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.75"
$ws.Range("E50").Value = "  +3.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.23"
$ws.Range("E51").Value = "  -1.28%  "
